# Corona Stats of Pakistan - add 2020-04-13 data block (rows 16-22)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new date block (2020-04-13), following the same CityNames /
# CityWiseCounts / Date / Headers / OverallCounts layout used by the
# existing 2020-04-09 and 2020-04-10 blocks.
$newRows = @(
    @{ Row = 16; A = "ICT";         B = "131";   C = "2020-04-13"; D = "Recovered";       E = "1,097" },
    @{ Row = 17; A = "Punjab";      B = "2,672"; C = "2020-04-13"; D = "Critical";         E = "44" },
    @{ Row = 18; A = "Sindh";       B = "1,452"; C = "2020-04-13"; D = "Deaths";           E = "93" },
    @{ Row = 19; A = "KP";          B = "744";   C = "2020-04-13"; D = "Cases (24 HRS)";   E = "122" },
    @{ Row = 20; A = "Balochistan"; B = "230";   C = "2020-04-13"; D = "Deaths (24 HRS)";  E = "7" },
    @{ Row = 21; A = "AJK";         B = "43";    C = "2020-04-13"; D = "Tests (24 HRS)";   E = "3,233" },
    @{ Row = 22; A = "GB";          B = "224";   C = "2020-04-13"; D = "Total Tests";      E = "65,114" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    foreach ($col in @("A", "B", "C", "D", "E")) {
        # Force text storage (matches the rest of the sheet, where numbers
        # and dates are stored as plain text, not numeric/date values).
        $ws.Range("${col}${rowNum}").NumberFormat = "@"
    }
    $ws.Range("A${rowNum}").Value = $r.A
    $ws.Range("B${rowNum}").Value = $r.B
    $ws.Range("C${rowNum}").Value = $r.C
    $ws.Range("D${rowNum}").Value = $r.D
    $ws.Range("E${rowNum}").Value = $r.E
    foreach ($col in @("A", "B", "C", "D", "E")) {
        # Drop the temporary text number-format again so no stray cell
        # style sticks around on the new cells.
        $ws.Range("${col}${rowNum}").ClearFormats()
    }
}

# Resize columns to fit the (now larger) content, mirroring the bestFit
# column widths Excel would compute for the refreshed data.
$ws.Columns.Item(1).ColumnWidth = 9.09
$ws.Columns.Item(2).ColumnWidth = 12.5
$ws.Columns.Item(3).ColumnWidth = 9.33
$ws.Columns.Item(4).ColumnWidth = 12.09
$ws.Columns.Item(5).ColumnWidth = 11.09

# Select the whole sheet (matches the saved selection state in the workbook).
$ws.Cells.Select()
